$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Natmi LR-pair data (Bmp6-Bmpr1a), following Dr Hou advice.
# Replaces the old 2-sender-cluster table (FAPs/sCs x 4 targets = 8 rows)
# with a new 3-sender-cluster table (ECs/FAPs/sCs x 4 targets = 12 rows).
$rows = @(
    @("ECs", "Bmp6", "Bmpr1a", "ECs", 2, 0.6666666666666666, 19.72450833333333, 59.173525, 0.5834853563809828, 0.5834853563809829, 3, 1, 3.767160666666667, 11.301482, 0.05327117368020191, 0.05327117368020191, 74.30539196267222, 668.74852766405, 0.03108294975962584, 0.03108294975962585),
    @("ECs", "Bmp6", "Bmpr1a", "FAPs", 2, 0.6666666666666666, 19.72450833333333, 59.173525, 0.5834853563809828, 0.5834853563809829, 3, 1, 39.446935, 118.340805, 0.5578165391591924, 0.5578165391591923, 778.0713981319584, 7002.642583187625, 0.3254777821465078, 0.3254777821465078),
    @("ECs", "Bmp6", "Bmpr1a", "M2", 2, 0.6666666666666666, 19.72450833333333, 59.173525, 0.5834853563809828, 0.5834853563809829, 3, 1, 1.169082333333333, 3.507247, 0.01653191714824366, 0.01653191714824366, 23.05957422618611, 207.536168035675, 0.009646131568903834, 0.009646131568903836),
    @("ECs", "Bmp6", "Bmpr1a", "sCs", 2, 0.6666666666666666, 19.72450833333333, 59.173525, 0.5834853563809828, 0.5834853563809829, 3, 1, 26.33350433333333, 79.000513, 0.372380370012362, 0.372380370012362, 519.4154256687028, 4674.738831018325, 0.2172784929059453, 0.2172784929059453),
    @("FAPs", "Bmp6", "Bmpr1a", "ECs", 3, 1, 0.7684289999999999, 2.305287, 0.02273146997336134, 0.02273146997336134, 3, 1, 3.767160666666667, 11.301482, 0.05327117368020191, 0.05327117368020191, 2.894795503926, 26.053159535334, 0.001210932084957227, 0.001210932084957227),
    @("FAPs", "Bmp6", "Bmpr1a", "FAPs", 3, 1, 0.7684289999999999, 2.305287, 0.02273146997336134, 0.02273146997336134, 3, 1, 39.446935, 118.340805, 0.5578165391591924, 0.5578165391591923, 30.312168815115, 272.809519336035, 0.01267998991054152, 0.01267998991054152),
    @("FAPs", "Bmp6", "Bmpr1a", "M2", 3, 1, 0.7684289999999999, 2.305287, 0.02273146997336134, 0.02273146997336134, 3, 1, 1.169082333333333, 3.507247, 0.01653191714824366, 0.01653191714824366, 0.8983567683209998, 8.085210914888998, 0.0003757947782573983, 0.0003757947782573983),
    @("FAPs", "Bmp6", "Bmpr1a", "sCs", 3, 1, 0.7684289999999999, 2.305287, 0.02273146997336134, 0.02273146997336134, 3, 1, 26.33350433333333, 79.000513, 0.372380370012362, 0.372380370012362, 20.235428401359, 182.118855612231, 0.008464753199605192, 0.008464753199605194),
    @("sCs", "Bmp6", "Bmpr1a", "ECs", 3, 1, 13.31169566666667, 39.935087, 0.3937831736456558, 0.3937831736456558, 3, 1, 3.767160666666667, 11.301482, 0.05327117368020191, 0.05327117368020191, 50.14729632210378, 451.3256668989341, 0.02097729183561884, 0.02097729183561884),
    @("sCs", "Bmp6", "Bmpr1a", "FAPs", 3, 1, 13.31169566666667, 39.935087, 0.3937831736456558, 0.3937831736456558, 3, 1, 39.446935, 118.340805, 0.5578165391591924, 0.5578165391591923, 525.1055937027817, 4725.950343325036, 0.219658767102143, 0.219658767102143),
    @("sCs", "Bmp6", "Bmpr1a", "M2", 3, 1, 13.31169566666667, 39.935087, 0.3937831736456558, 0.3937831736456558, 3, 1, 1.169082333333333, 3.507247, 0.01653191714824366, 0.01653191714824366, 15.56246823060989, 140.062214075489, 0.00650999080108243, 0.006509990801082431),
    @("sCs", "Bmp6", "Bmpr1a", "sCs", 3, 1, 13.31169566666667, 39.935087, 0.3937831736456558, 0.3937831736456558, 3, 1, 26.33350433333333, 79.000513, 0.372380370012362, 0.372380370012362, 350.5435955221812, 3154.892359699631, 0.1466371239068115, 0.1466371239068115)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $rowData = $rows[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}

Write-Host "Updated Bmp6-Bmpr1a LR-pair rows 2-13 with new NATMI output"